$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.45491533333333
$ws.Range("H2").Value = 31.364746
$ws.Range("I2").Value = 0.0134573334963438
$ws.Range("J2").Value = 0.0134573334963438
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.4549469999999
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 3113.66527501094
$ws.Range("R2").Value = 28022.98747509846
$ws.Range("S2").Value = 0.01191264210127866
$ws.Range("T2").Value = 0.01191264210127866

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.45491533333333
$ws.Range("H3").Value = 31.364746
$ws.Range("I3").Value = 0.0134573334963438
$ws.Range("J3").Value = 0.0134573334963438
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 254.4762391891798
$ws.Range("R3").Value = 2290.286152702618
$ws.Range("S3").Value = 0.0009736063747987267
$ws.Range("T3").Value = 0.0009736063747987267

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.45491533333333
$ws.Range("H4").Value = 31.364746
$ws.Range("I4").Value = 0.0134573334963438
$ws.Range("J4").Value = 0.0134573334963438
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 149.267272664188
$ws.Range("R4").Value = 1343.405453977692
$ws.Range("S4").Value = 0.0005710850202664125
$ws.Range("T4").Value = 0.0005710850202664125

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 735.4993083333334
$ws.Range("H5").Value = 2206.497925
$ws.Range("I5").Value = 0.9467182815928301
$ws.Range("J5").Value = 0.9467182815928301
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.4549469999999
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 219045.1651818317
$ws.Range("R5").Value = 1971406.486636485
$ws.Range("S5").Value = 0.8380498307794045
$ws.Range("T5").Value = 0.8380498307794046

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 735.4993083333334
$ws.Range("H6").Value = 2206.497925
$ws.Range("I6").Value = 0.9467182815928301
$ws.Range("J6").Value = 0.9467182815928301
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 17902.30642176184
$ws.Range("R6").Value = 161120.7577958565
$ws.Range("S6").Value = 0.06849283733272263
$ws.Range("T6").Value = 0.06849283733272264

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 735.4993083333334
$ws.Range("H7").Value = 2206.497925
$ws.Range("I7").Value = 0.9467182815928301
$ws.Range("J7").Value = 0.9467182815928301
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 10500.89573191315
$ws.Range("R7").Value = 94508.06158721836
$ws.Range("S7").Value = 0.04017561348070289
$ws.Range("T7").Value = 0.04017561348070289

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 30.939307
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.03982438491082609
$ws.Range("J8").Value = 0.03982438491082609
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.4549469999999
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 9214.292298633911
$ws.Range("R8").Value = 82928.6306877052
$ws.Range("S8").Value = 0.03525316842858402
$ws.Range("T8").Value = 0.03525316842858402

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 30.939307
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.03982438491082609
$ws.Range("J9").Value = 0.03982438491082609
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 753.0733858147105
$ws.Range("R9").Value = 6777.660472332394
$ws.Range("S9").Value = 0.002881200427421431
$ws.Range("T9").Value = 0.002881200427421431

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 30.939307
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.03982438491082609
$ws.Range("J10").Value = 0.03982438491082609
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 441.727725836838
$ws.Range("R10").Value = 3975.549532531542
$ws.Range("S10").Value = 0.001690016054820634
$ws.Range("T10").Value = 0.001690016054820635
